# Updated cryptos list on Tue Nov 21 03:14:49 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be stored as text, matching the source data,
# since Excel would otherwise auto-convert plain-looking numbers (e.g. "55.96")
# into numeric values.
$priceCells = "D2","D3","D5","D6","D8","D9","D10","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D38","D39","D41","D43","D44","D45","D46","D47","D48","D49","D50"
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "37.540.82"
$ws.Range("E2").Value = "  +0.59%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.014.51"
$ws.Range("E3").Value = "  +0.36%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.14%  "

# Row 5 - BNB
$ws.Range("D5").Value = "262.88"
$ws.Range("E5").Value = "  +6.22%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.617"
$ws.Range("E6").Value = "  -1.25%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.10%  "

# Row 8 - Solana
$ws.Range("D8").Value = "55.96"
$ws.Range("E8").Value = "  -7.40%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.385"
$ws.Range("E9").Value = "  +0.15%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.0774"
$ws.Range("E10").Value = "  -3.73%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -2.13%  "

# Row 12 - was WrappedliquidstakedEther2.0, now Chainlink
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "14.33"
$ws.Range("E12").Value = "  -5.05%  "

# Row 13 - was Chainlink, now WrappedliquidstakedEther2.0
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.308.63"
$ws.Range("E13").Value = "  +0.23%  "

# Row 14 - Polygon
$ws.Range("D14").Value = "0.805"
$ws.Range("E14").Value = "  -5.26%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "20.80"
$ws.Range("E15").Value = "  -8.87%  "

# Row 16 - Polkadot
$ws.Range("D16").Value = "5.26"
$ws.Range("E16").Value = "  -3.97%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "1.993.21"
$ws.Range("E17").Value = "  -1.05%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "37.474.76"
$ws.Range("E18").Value = "  +0.73%  "

# Row 19 - Litecoin
$ws.Range("D19").Value = "69.43"
$ws.Range("E19").Value = "  -1.64%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0842"
$ws.Range("E20").Value = "  -3.07%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "5.17"
$ws.Range("E21").Value = "  -0.96%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "228.22"
$ws.Range("E22").Value = "  -1.27%  "

# Row 23 - PancakeSwap
$ws.Range("D23").Value = "2.69"
$ws.Range("E23").Value = "  +7.65%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.05%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "2.33"
$ws.Range("E25").Value = "  -1.49%  "

# Row 26 - Monero
$ws.Range("D26").Value = "163.65"
$ws.Range("E26").Value = "  -0.08%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "8.95"
$ws.Range("E27").Value = "  -5.33%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "19.66"
$ws.Range("E28").Value = "  -0.40%  "

# Row 29 - Kaspa
$ws.Range("D29").Value = "0.127"
$ws.Range("E29").Value = "  -12.38%  "

# Row 30 - ImmutableX
$ws.Range("D30").Value = "1.34"
$ws.Range("E30").Value = "  -0.71%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "0.120"
$ws.Range("E31").Value = "  -1.43%  "

# Row 32 - Hedera
$ws.Range("D32").Value = "0.0651"
$ws.Range("E32").Value = "  -1.24%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "4.63"
$ws.Range("E33").Value = "  -4.49%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").Value = "4.53"
$ws.Range("E34").Value = "  -0.26%  "

# Row 35 - LidoDAOToken
$ws.Range("D35").Value = "2.40"
$ws.Range("E35").Value = "  +0.89%  "

# Row 36 - WEMIXToken
$ws.Range("E36").Value = "  +1.20%  "

# Row 37 - BinanceUSD
$ws.Range("E37").Value = "  -0.14%  "

# Row 38 - RenderToken
$ws.Range("D38").Value = "3.35"
$ws.Range("E38").Value = "  +2.14%  "

# Row 39 - THORChain
$ws.Range("D39").Value = "5.20"
$ws.Range("E39").Value = "  -5.22%  "

# Row 40 - HuobiToken
$ws.Range("E40").Value = "  +4.56%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "1.22"
$ws.Range("E41").Value = "  +2.65%  "

# Row 42 - Cronos
$ws.Range("E42").Value = "  -4.44%  "

# Row 43 - VeChain
$ws.Range("D43").Value = "0.0213"
$ws.Range("E43").Value = "  -1.10%  "

# Row 44 - Maker
$ws.Range("D44").Value = "1.398.10"

# Row 45 - Aave
$ws.Range("D45").Value = "89.87"
$ws.Range("E45").Value = "  -1.33%  "

# Row 46 - InjectiveProtocol
$ws.Range("D46").Value = "15.65"
$ws.Range("E46").Value = "  -6.54%  "

# Row 47 - ARBITRUM
$ws.Range("D47").Value = "1.03"
$ws.Range("E47").Value = "  -2.26%  "

# Row 48 - FraxShare
$ws.Range("D48").Value = "7.07"
$ws.Range("E48").Value = "  -3.02%  "

# Row 49 - MXToken
$ws.Range("D49").Value = "2.88"
$ws.Range("E49").Value = "  +0.87%  "

# Row 50 - RocketPoolETH
$ws.Range("D50").Value = "2.199.98"
$ws.Range("E50").Value = "  +0.14%  "

# Row 51 - NEARProtocol
$ws.Range("E51").Value = "  -3.49%  "
